# Insert 6 new weekly-report rows (new price-report date 44449) right before
# the existing block that starts at row 795, shifting the rest of the data
# down by 6 rows (795-879 -> 801-885) and extending the used range to R885.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at 795..800 (existing rows shift down to 801..885).
$ws.Range("A795:A800").EntireRow.Insert()

# Static / repeated metadata columns, identical for all six new rows.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112037
$categoria = "Cebollín"
$variedad  = "Sin especificar"
$unidad    = "`$/paquete 36 unidades"
$kgUnid    = 36
$clasif    = "Hortaliza"
$fecha     = 44449

# Per-row data: row number, Calidad, Volumen, Precio min, Precio max, Precio prom, Origen, Precio $/Kg
$rows = @(
    @{ R = 795; I = "Extra";   J = 670; K = 3000; L = 3300; M = 3166; O = "Provincia de Chacabuco"; P = 88 },
    @{ R = 796; I = "Extra";   J = 680; K = 3000; L = 3300; M = 3172; O = "Región Metropolitana";   P = 88 },
    @{ R = 797; I = "Primera"; J = 850; K = 2400; L = 2600; M = 2506; O = "Provincia de Chacabuco"; P = 70 },
    @{ R = 798; I = "Primera"; J = 790; K = 2400; L = 2600; M = 2511; O = "Región Metropolitana";   P = 70 },
    @{ R = 799; I = "Segunda"; J = 250; K = 1800; L = 1800; M = 1800; O = "Provincia de Chacabuco"; P = 50 },
    @{ R = 800; I = "Segunda"; J = 310; K = 1900; L = 1900; M = 1900; O = "Región Metropolitana";   P = 53 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $catId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $variedad
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $kgUnid
    $ws.Cells.Item($r, 18).Value = $clasif
}
